$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (ID "H 72"); this shifts all subsequent rows
# up by one and shrinks the used range from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()

# After the shift, a handful of cells swap between having a value and
# being empty (missing-data pattern changed). Apply those adjustments.

# Column B ("A" in the source header)
$ws.Range("B9").Value = -19.9
$ws.Range("B10").ClearContents()
$ws.Range("B11").Value = -19.9
$ws.Range("B12").ClearContents()
$ws.Range("B17").Value = -19.9
$ws.Range("B18").ClearContents()
$ws.Range("B36").Value = -19.1
$ws.Range("B37").ClearContents()
$ws.Range("B54").Value = -17.2
$ws.Range("B55").ClearContents()

# Column E ("D" in the source header)
$ws.Range("E5").Value = -5.7
$ws.Range("E7").ClearContents()
$ws.Range("E14").Value = -7.9
$ws.Range("E15").Value = -12
$ws.Range("E16").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("E45").Value = -7.4
$ws.Range("E47").ClearContents()
$ws.Range("E53").Value = -5.7
$ws.Range("E55").ClearContents()
$ws.Range("E58").Value = -6.8
$ws.Range("E60").ClearContents()
